$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Add the new test-case row (row 5).
# NB: the cells are populated in this particular order so that the three
# brand-new shared-string entries get appended to the shared string table
# in the same order as the target workbook (OPQA-3455, description, RCC111).
$ws.Range("B5").Value = "OPQA-3455"
$ws.Range("C5").Value = "Verify that user is able to add an article to the group from search results page."
$ws.Range("A5").Value = "RCC111"
$ws.Range("D5").Value = "Y"

# Rows 3 and 4 (already present) and the new row 5 should now carry the same
# border formatting that row 2 already has.
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A3:E5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update the active selection shown in the sheet view.
$ws.Range("C18").Select() | Out-Null
